# Updates cryptos list (row data) to reflect the latest price/volume snapshot.
# Values in column D that look like plain numbers are written with a leading
# apostrophe so Excel keeps them as text (preserving formatting such as
# trailing zeros or thousands-style separators), matching the source data
# which stores Price/Volume as text strings, not numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '49.709.58'
$ws.Range("E2").Value = '  -0.68%  '

$ws.Range("D3").Value = '2.641.14'
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''113.50'
$ws.Range("E5").Value = '  +2.36%  '

$ws.Range("D6").Value = '''324.61'
$ws.Range("E6").Value = '  -0.93%  '

$ws.Range("E7").Value = '  -1.10%  '

$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").Value = '''0.546'
$ws.Range("E9").Value = '  -2.75%  '

$ws.Range("D10").Value = '''40.01'
$ws.Range("E10").Value = '  -1.88%  '

$ws.Range("D11").Value = '''19.84'
$ws.Range("E11").Value = '  -4.06%  '

$ws.Range("E12").Value = '  -0.90%  '

$ws.Range("E13").Value = '  +1.10%  '

$ws.Range("E14").Value = '  +0.40%  '

$ws.Range("D15").Value = '3.054.59'
$ws.Range("E15").Value = '  +0.06%  '

$ws.Range("D16").Value = '2.650.49'
$ws.Range("E16").Value = '  +0.37%  '

$ws.Range("E17").Value = '  -2.20%  '

$ws.Range("D18").Value = '49.608.80'
$ws.Range("E18").Value = '  -0.86%  '

$ws.Range("D19").Value = '''2.98'
$ws.Range("E19").Value = '  -2.90%  '

$ws.Range("D20").Value = '''12.93'
$ws.Range("E20").Value = '  -2.93%  '

$ws.Range("E21").Value = '  -1.62%  '

$ws.Range("D22").Value = '0.0₃0948'
$ws.Range("E22").Value = '  -1.54%  '

$ws.Range("D23").Value = '''270.23'
$ws.Range("E23").Value = '  -3.53%  '

$ws.Range("D24").Value = '''69.00'
$ws.Range("E24").Value = '  -5.52%  '

$ws.Range("D25").Value = '''2.55'
$ws.Range("E25").Value = '  -2.00%  '

$ws.Range("D26").Value = '''26.40'
$ws.Range("E26").Value = '  -0.91%  '

$ws.Range("E27").Value = '  +0.06%  '

$ws.Range("D28").Value = '''10.36'
$ws.Range("E28").Value = '  +4.11%  '

$ws.Range("D29").Value = '''2.22'
$ws.Range("E29").Value = '  -0.47%  '

$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").Value = '''0.140'
$ws.Range("E30").Value = '  -3.51%  '

$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = '''35.28'
$ws.Range("E31").Value = '  -3.60%  '

$ws.Range("D32").Value = '''49.63'
$ws.Range("E32").Value = '  -0.19%  '

$ws.Range("D33").Value = '''5.49'
$ws.Range("E33").Value = '  +0.77%  '

$ws.Range("D34").Value = '''0.0815'
$ws.Range("E34").Value = '  +2.31%  '

$ws.Range("E35").Value = '  -0.33%  '

$ws.Range("D36").Value = '''19.08'
$ws.Range("E36").Value = '  -3.66%  '

$ws.Range("E37").Value = '  +4.18%  '

$ws.Range("D38").Value = '''2.06'
$ws.Range("E38").Value = '  -0.47%  '

$ws.Range("E39").Value = '  +0.94%  '

$ws.Range("D40").Value = '''126.99'
$ws.Range("E40").Value = '  +2.84%  '

$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '''22.49'
$ws.Range("E41").Value = '  +0.14%  '

$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").Value = '''0.111'
$ws.Range("E42").Value = '  -1.40%  '

$ws.Range("E43").Value = '  +4.28%  '

$ws.Range("E44").Value = '  -3.22%  '

$ws.Range("D45").Value = '2.063.05'
$ws.Range("E45").Value = '  -0.13%  '

$ws.Range("D46").Value = '''3.23'
$ws.Range("E46").Value = '  -3.26%  '

$ws.Range("D47").Value = '''2.14'
$ws.Range("E47").Value = '  +6.64%  '

$ws.Range("E48").Value = '  -7.00%  '

$ws.Range("D49").Value = '''8.95'
$ws.Range("E49").Value = '  -1.36%  '

$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").Value = '''5.23'
$ws.Range("E50").Value = '  -3.03%  '

$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = '''59.25'
$ws.Range("E51").Value = '  +1.41%  '

